$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("summary")
$ws2 = $wb.Worksheets.Item("model_fit")
$ws3 = $wb.Worksheets.Item("steps")

# ---- sheet3 "steps": update text (shared-string) cells ----
# Ensure these cells keep Text format so numeric-looking strings
# (e.g. "-0.11", "-0.30") are not auto-converted to numbers.
$stepsCells = @("B2","C2","D2","B3","C3","B4","C4","B5","C5","D5")
foreach ($addr in $stepsCells) {
    $ws3.Range($addr).NumberFormat = "@"
}

$ws3.Range("B2").Value = "1.67 (0.072)"
$ws3.Range("C2").Value = "-1.56 (0.072)"
$ws3.Range("D2").Value = "-0.11"
$ws3.Range("B3").Value = "0.36 (0.072)"
$ws3.Range("C3").Value = "-0.36"
$ws3.Range("B4").Value = "0.84 (0.082)"
$ws3.Range("C4").Value = "-0.84"
$ws3.Range("B5").Value = "1.26 (0.082)"
$ws3.Range("C5").Value = "-0.96 (0.092)"
$ws3.Range("D5").Value = "-0.30"

# ---- sheet1 "summary": numeric updates ----
$ws1.Range("D2").Value = 1440
$ws1.Range("E2").Value = 82.22
$ws1.Range("F2").Value = -1.96
$ws1.Range("H2").Value = 1.05
$ws1.Range("I2").Value = 1.03
$ws1.Range("J2").Value = 0.27
$ws1.Range("L2").Value = 1
$ws1.Range("D3").Value = 1431
$ws1.Range("E3").Value = 79.11
$ws1.Range("I3").Value = 1.03
$ws1.Range("K3").Value = 0.04
$ws1.Range("D4").Value = 1430
$ws1.Range("F4").Value = -1.56
$ws1.Range("I4").Value = -1.68
$ws1.Range("L4").Value = 1.29
$ws1.Range("D5").Value = 1428
$ws1.Range("E5").Value = 70.66
$ws1.Range("F5").Value = -1.14
$ws1.Range("H5").Value = 1.03
$ws1.Range("I5").Value = 0.83
$ws1.Range("D6").Value = 1429
$ws1.Range("E6").Value = 67.32
$ws1.Range("F6").Value = -0.93
$ws1.Range("I6").Value = 1.37
$ws1.Range("J6").Value = 0.35
$ws1.Range("L6").Value = 1.11
$ws1.Range("D7").Value = 1425
$ws1.Range("E7").Value = 63.79
$ws1.Range("F7").Value = -0.73
$ws1.Range("H7").Value = 0.99
$ws1.Range("I7").Value = -0.24
$ws1.Range("L7").Value = 1.31
$ws1.Range("D8").Value = 1415
$ws1.Range("H8").Value = 0.91
$ws1.Range("I8").Value = -2.46
$ws1.Range("J8").Value = 0.43
$ws1.Range("K8").Value = 0.04
$ws1.Range("L8").Value = 0.98
$ws1.Range("D9").Value = 1397
$ws1.Range("E9").Value = 53.26
$ws1.Range("I9").Value = 0.96
$ws1.Range("J9").Value = 0.36
$ws1.Range("D10").Value = 1388
$ws1.Range("E10").Value = 52.45
$ws1.Range("F10").Value = -0.12
$ws1.Range("H10").Value = 0.99
$ws1.Range("I10").Value = -0.43
$ws1.Range("J10").Value = 0.42
$ws1.Range("K10").Value = 0.05
$ws1.Range("L10").Value = 1.35
$ws1.Range("D11").Value = 1349
$ws1.Range("E11").Value = 45.74
$ws1.Range("F11").Value = 0.24
$ws1.Range("H11").Value = 1.01
$ws1.Range("I11").Value = 0.29
$ws1.Range("J11").Value = 0.39
$ws1.Range("L11").Value = 1.23
$ws1.Range("D12").Value = 1308
$ws1.Range("E12").Value = 39.45
$ws1.Range("F12").Value = 0.55
$ws1.Range("H12").Value = 1.02
$ws1.Range("I12").Value = 0.74
$ws1.Range("J12").Value = 0.35
$ws1.Range("K12").Value = 0.03
$ws1.Range("L12").Value = 1.15
$ws1.Range("D13").Value = 1246
$ws1.Range("E13").Value = 36.84
$ws1.Range("F13").Value = 0.71
$ws1.Range("H13").Value = 1.02
$ws1.Range("I13").Value = 0.58
$ws1.Range("J13").Value = 0.36
$ws1.Range("K13").Value = 0.04
$ws1.Range("L13").Value = 1.18
$ws1.Range("D14").Value = 1178
$ws1.Range("E14").Value = 34.04
$ws1.Range("F14").Value = 0.86
$ws1.Range("H14").Value = 1.06
$ws1.Range("I14").Value = 1.81
$ws1.Range("J14").Value = 0.32
$ws1.Range("L14").Value = 0.99
$ws1.Range("D15").Value = 1085
$ws1.Range("F15").Value = 0.17
$ws1.Range("J15").Value = 0.46
$ws1.Range("L15").Value = 0.79
$ws1.Range("D16").Value = 946
$ws1.Range("E16").Value = 22.3
$ws1.Range("F16").Value = 1.61
$ws1.Range("G16").Value = 0.09
$ws1.Range("H16").Value = 1.08
$ws1.Range("I16").Value = 1.74
$ws1.Range("J16").Value = 0.29
$ws1.Range("L16").Value = 0.92
$ws1.Range("D17").Value = 763
$ws1.Range("F17").Value = 0.24
$ws1.Range("H17").Value = 0.93
$ws1.Range("I17").Value = -1.73
$ws1.Range("J17").Value = 0.54
$ws1.Range("K17").Value = 0.06
$ws1.Range("L17").Value = 0.82
$ws1.Range("D18").Value = 460
$ws1.Range("E18").Value = 21.52
$ws1.Range("F18").Value = 1.7
$ws1.Range("G18").Value = 0.13
$ws1.Range("H18").Value = 0.97
$ws1.Range("I18").Value = -0.36
$ws1.Range("J18").Value = 0.38
$ws1.Range("K18").Value = 0.04
$ws1.Range("L18").Value = 1.43

# ---- sheet2 "model_fit": numeric updates ----
$ws2.Range("D2").Value = 26064
$ws2.Range("E2").Value = 26112
$ws2.Range("F2").Value = 26240
$ws2.Range("G2").Value = 0.766
$ws2.Range("H2").Value = 0.701
$ws2.Range("D3").Value = 25941
$ws2.Range("E3").Value = 26021
$ws2.Range("F3").Value = 26233
$ws2.Range("G3").Value = 0.774
$ws2.Range("H3").Value = 0.704

Write-Host "Edits applied"
